$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the end date for the "Git Depostory erstellen" task (column G, row 1)
$ws.Range("G1").Value = 43151

# Widen column A to fit the new, longer task description
$ws.Columns.Item(1).ColumnWidth = 36.5

# Mark "S" for the task in row 7 across columns D and E
$ws.Range("D7").Value = "S"
$ws.Range("E7").Value = "S"

# Add a new task row (12) for the use case diagram
$ws.Range("A12").Value = "Erstellen eines Anwendungsfalldiagramm"
$ws.Range("G12").Value = "S"

# Update the active selection to reflect where the user left off
$ws.Activate()
$ws.Range("F15").Select()
